$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 3, pushing existing rows 3-16 down to 4-17.
$ws.Rows("3:3").Insert()

# 2) Copy formatting from the row below (old row3, now row4) into new row3.
$ws.Range("A4:I4").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122) # xlPasteFormats

# 3) Populate new row 3 with the gpt-4o data.
$ws.Range("A3").Value = "gpt-4o"
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = 0.422
$ws.Range("D3").Value = 0.3548
$ws.Range("E3").Value = 0.3492
$ws.Range("F3").Value = 0.9411
$ws.Range("G3").Value = 0.941
$ws.Range("H3").Value = 0.9411
$ws.Range("I3").Formula = "=AVERAGE(H3,E3)"

# 4) Fill helper column O with tracked values (rows 1-68).
$vals = @(1,1,1,1,1,1,1,1,1,1,0,1,1,1,0,1,1,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,1,1,1,1)
$n = $vals.Length
$arr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
  $arr[$i,0] = $vals[$i]
}
$ws.Range("O1:O68").Value = $arr

Write-Host "done"
